$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Total" column header (N1) ---
$ws.Range("J1").Copy()
$ws.Range("N1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("N1").Value = "Total"

# --- Total formula for existing row 2 ---
$ws.Range("C2").Copy()
$ws.Range("N2").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("N2").Formula = "=C2+E2+G2+I2+K2+M2"

# --- New row 3: copy formatting from row 2 for columns B:M, N ---
$ws.Range("B2:M2").Copy()
$ws.Range("B3:M3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("N2").Copy()
$ws.Range("N3").PasteSpecial(-4122) # xlPasteFormats

# A3 date (keep the underlined font already on A3, add the date number format)
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("A3").Value = 45083.338888888888

$ws.Range("B3").Value = 0.38
$ws.Range("C3").Formula = "=B3 * 600 / 100 *100"

$ws.Range("D3").Value = 0.92
$ws.Range("E3").Formula = "=D3 * 600 / 100 * 100"

$ws.Range("F3").Value = 0.51
$ws.Range("G3").Formula = "=F3 * 600 / 100 * 100"

$ws.Range("H3").Value = 0.27
$ws.Range("I3").Formula = "=H3 * 600 / 100 * 100"

$ws.Range("J3").Value = 0.84
$ws.Range("K3").Formula = "=J3 * 600 / 100 * 100"

$ws.Range("L3").Value = 0.18
$ws.Range("M3").Formula = "=L3 * 600 / 100 * 100"

$ws.Range("N3").Formula = "=C3+E3+G3+I3+K3+M3"

# --- Column width for new column N (stored width=11 => ColumnWidth = 11 - 5/6) ---
$ws.Columns("N").ColumnWidth = 10.166666666666666

# --- Selection as recorded after the edit ---
$ws.Range("P3").Select()
